# Weekly update: insert a new daily price record as row 4, pushing the
# existing data rows down by one (row 4 -> row 5, row 5 -> row 6, ... row 99 -> row 100).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4 (shifts rows 4..99 down to 5..100).
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the new record's data.
$ws.Cells.Item(4, 1).Value  = 10
$ws.Cells.Item(4, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(4, 3).Value  = "La Araucanía"
$ws.Cells.Item(4, 4).Value  = 44515
$ws.Cells.Item(4, 5).Value  = 9
$ws.Cells.Item(4, 6).Value  = 100112012
$ws.Cells.Item(4, 7).Value  = "Espinaca"
$ws.Cells.Item(4, 8).Value  = "Sin especificar"
$ws.Cells.Item(4, 9).Value  = "Primera"
$ws.Cells.Item(4, 10).Value = 45
$ws.Cells.Item(4, 11).Value = 8000
$ws.Cells.Item(4, 12).Value = 8000
$ws.Cells.Item(4, 13).Value = 8000
$ws.Cells.Item(4, 14).Value = '$/docena de atados'
$ws.Cells.Item(4, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(4, 16).Value = 2667
$ws.Cells.Item(4, 17).Value = 3
$ws.Cells.Item(4, 18).Value = "Hortaliza"

# Give the new date cell the same date number format as the rest of column D.
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
